{"js": "// The document's \"content slots\" get shuffled: several paragraphs (and a\n// few runs inside the \"Avalia\u00e7\u00e3o\" bullet paragraph) keep their position,\n// style and run-level formatting (bold labels, italic EN translations,\n// list styles, etc.) but their *text* moves to a different slot.\n//\n// Capture every source text first (before any paragraph is touched) and\n// only then write the new values back, so a replacement can never clobber\n// a value we still need to read.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Paragraph-level single-run text slots (index -> Range to replace).\nconst idx = [5, 6, 8, 10, 11, 13, 18];\nconst ranges = {};\nfor (const i of idx) {\n  ranges[i] = paragraphs.items[i].getRange();\n}\nfor (const i of idx) {\n  ranges[i].load(\"text\");\n}\n\n// The \"Avalia\u00e7\u00e3o\" bullet paragraph (index 16) holds three label/value run\n// pairs (\"M\u00e9todo: \"/\"Crit\u00e9rio: \"/\"Norma de recupera\u00e7\u00e3o: \"); only the value\n// runs move. Locate them unambiguously via search on their known text.\nconst metodoValue = body.search(\"Aulas expositivas, semin\u00e1rios e exerc\u00edcios comentados.\", { matchCase: true });\nconst criterioValue = body.search(\"M\u00e9dia aritm\u00e9tica de duas provas sendo a primeira com peso 1 e a segunda com peso 2.\", { matchCase: true });\nconst normaValue = body.search(\"Aplica\u00e7\u00e3o de uma prova escrita dentro do prazo regimental antes do in\u00edcio do pr\u00f3ximo semestre letivo. A nota da segunda avalia\u00e7\u00e3o ser\u00e1 a m\u00e9dia aritm\u00e9tica entre a nota da prova de recupera\u00e7\u00e3o e a nota final da primeira avalia\u00e7\u00e3o\", { matchCase: true });\nmetodoValue.load(\"text\");\ncriterioValue.load(\"text\");\nnormaValue.load(\"text\");\n\nawait context.sync();\n\nconst text = {};\nfor (const i of idx) {\n  text[i] = ranges[i].text;\n}\nconst metodoText = metodoValue.items[0].text;\nconst criterioText = criterioValue.items[0].text;\nconst normaText = normaValue.items[0].text;\n\n// New values, per the target layout (old slot -> new slot):\n//  para5  \"Objetivos\" PT text              -> para8  Docente bullet\n//  para6  \"Objetivos\" EN text (italic)     -> para11 \"Programa resumido\" EN text\n//  para8  Docente bullet (\"144651...\")     -> para18 Bibliografia paragraph\n//  para10 \"Programa resumido\" PT text      -> para5  \"Objetivos\" PT text\n//  para11 \"Programa resumido\" EN text      -> para6  \"Objetivos\" EN text (italic)\n//  para13 \"Programa\" PT text               -> para10 \"Programa resumido\" PT text\n//  para16 M\u00e9todo value (\"Aulas...\")        -> para13 \"Programa\" PT text\n//  para16 Crit\u00e9rio value (\"M\u00e9dia...\")      -> para16 M\u00e9todo value\n//  para16 Norma value (\"Aplica\u00e7\u00e3o...\")     -> para16 Crit\u00e9rio value\n//  para18 Bibliografia paragraph (BOCKRIS) -> para16 Norma value\nranges[5].insertText(text[10], \"Replace\");\nranges[6].insertText(text[11], \"Replace\");\nranges[8].insertText(text[5], \"Replace\");\nranges[10].insertText(text[13], \"Replace\");\nranges[11].insertText(text[6], \"Replace\");\nranges[13].insertText(metodoText, \"Replace\");\nranges[18].insertText(text[8], \"Replace\");\n\nmetodoValue.items[0].insertText(criterioText, \"Replace\");\ncriterioValue.items[0].insertText(normaText, \"Replace\");\nnormaValue.items[0].insertText(text[18], \"Replace\");\n\nawait context.sync();\n", "ps1": "# The document's \"content slots\" get shuffled: several paragraphs (and a\n# few runs inside the \"Avalia\u00e7\u00e3o\" bullet paragraph) keep their position,\n# style and run-level formatting (bold labels, italic EN translations,\n# list styles, etc.) but their *text* moves to a different slot.\n#\n# Strategy: capture every source text first (before any paragraph is\n# touched), stamp each source location with a unique placeholder token so\n# later Find/Replace passes can never match the wrong occurrence, and only\n# then swap the placeholders for their final destination text.\n\n$d = $word.ActiveDocument\n\n# Paragraphs.Item is 1-based; these are the 0-based indices (matching the\n# Office.js / OOXML paragraph order) of the single-run paragraphs whose\n# whole text moves elsewhere.\n$p5  = $d.Paragraphs.Item(6)    # \"Objetivos\" PT text\n$p6  = $d.Paragraphs.Item(7)    # \"Objetivos\" EN text (italic)\n$p8  = $d.Paragraphs.Item(9)    # Docente bullet (\"144651 - Antonio Fernando Sartori\")\n$p10 = $d.Paragraphs.Item(11)   # \"Programa resumido\" PT text\n$p11 = $d.Paragraphs.Item(12)   # \"Programa resumido\" EN text (italic)\n$p13 = $d.Paragraphs.Item(14)   # \"Programa\" PT text\n$p18 = $d.Paragraphs.Item(19)   # Bibliografia paragraph\n\nfunction Get-ParaText($para) {\n  $t = $para.Range.Text\n  return $t.Substring(0, $t.Length - 1)   # drop the trailing paragraph mark\n}\n\n$text5  = Get-ParaText $p5\n$text6  = Get-ParaText $p6\n$text8  = Get-ParaText $p8\n$text10 = Get-ParaText $p10\n$text11 = Get-ParaText $p11\n$text13 = Get-ParaText $p13\n$text18 = Get-ParaText $p18\n\n# The \"Avalia\u00e7\u00e3o\" bullet paragraph holds three label/value run pairs\n# (\"M\u00e9todo: \"/\"Crit\u00e9rio: \"/\"Norma de recupera\u00e7\u00e3o: \"); only the value runs\n# move. Their text is known (and unique) ahead of time.\n$metodoText   = \"Aulas expositivas, semin\u00e1rios e exerc\u00edcios comentados.\"\n$criterioText = \"M\u00e9dia aritm\u00e9tica de duas provas sendo a primeira com peso 1 e a segunda com peso 2.\"\n$normaText    = \"Aplica\u00e7\u00e3o de uma prova escrita dentro do prazo regimental antes do in\u00edcio do pr\u00f3ximo semestre letivo. A nota da segunda avalia\u00e7\u00e3o ser\u00e1 a m\u00e9dia aritm\u00e9tica entre a nota da prova de recupera\u00e7\u00e3o e a nota final da primeira avalia\u00e7\u00e3o\"\n\nfunction Replace-DocText($findText, $replaceText) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $findText\n  $find.Replacement.Text = $replaceText\n  $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 0, $false, $find.Replacement.Text, 2) | Out-Null\n}\n\n# --- Phase 1: stamp every source location with a unique placeholder ----\nReplace-DocText $text5        \"@@SLOT_P5@@\"\nReplace-DocText $text6        \"@@SLOT_P6@@\"\nReplace-DocText $text8        \"@@SLOT_P8@@\"\nReplace-DocText $text10       \"@@SLOT_P10@@\"\nReplace-DocText $text11       \"@@SLOT_P11@@\"\nReplace-DocText $text13       \"@@SLOT_P13@@\"\nReplace-DocText $metodoText   \"@@SLOT_METODO@@\"\nReplace-DocText $criterioText \"@@SLOT_CRITERIO@@\"\nReplace-DocText $normaText    \"@@SLOT_NORMA@@\"\nReplace-DocText $text18       \"@@SLOT_P18@@\"\n\n# --- Phase 2: place each captured text into its final destination ------\n# old slot (placeholder)  -> final text that belongs there\nReplace-DocText \"@@SLOT_P5@@\"       $text10\nReplace-DocText \"@@SLOT_P6@@\"       $text11\nReplace-DocText \"@@SLOT_P8@@\"       $text5\nReplace-DocText \"@@SLOT_P10@@\"      $text13\nReplace-DocText \"@@SLOT_P11@@\"      $text6\nReplace-DocText \"@@SLOT_P13@@\"      $metodoText\nReplace-DocText \"@@SLOT_P18@@\"      $text8\nReplace-DocText \"@@SLOT_METODO@@\"   $criterioText\nReplace-DocText \"@@SLOT_CRITERIO@@\" $normaText\nReplace-DocText \"@@SLOT_NORMA@@\"    $text18\n"}
